# Auto-generated edit script applying cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.738.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.983.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.94%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "497.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.33%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.981.68"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.97%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.425"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.31%  "

$ws.Range("E11").Value = "  -4.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.350"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.11%  "

$ws.Range("E13").Value = "  -0.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.491.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.702.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.982.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.96%  "

$ws.Range("E18").Value = "  -4.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.25%  "

$ws.Range("E21").Value = "  -2.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "324.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.85%  "

$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.463"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "61.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -10.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.162"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0904"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.61%  "

$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.81%  "

$ws.Range("E33").Value = "  -7.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.68%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.53%  "

$ws.Range("E37").Value = "  -6.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.17%  "

$ws.Range("E39").Value = "  -3.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.017.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("E44").Value = "  -6.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.636"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.202.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0236"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.87%  "
